# Weekly fruit/vegetable price update:
# Insert two new rows (296:297) for the latest reporting date, pushing the
# existing historical rows down by two. The previously-last pair of rows
# (old 416/417) end up duplicated at the very end as part of the shift,
# and rows 296/297 are populated with the newest observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 296:417 down to 298:419, creating two blank rows at 296:297.
$ws.Rows("296:297").Insert()

# --- Row 296 (Primera / $/paquete 4 unidades) ---
$ws.Range("A296").Value = 1
$ws.Range("B296").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C296").Value = "Arica y Parinacota"
$ws.Range("D296").Value = 44924
$ws.Range("E296").Value = 15
$ws.Range("F296").Value = 100114014
$ws.Range("G296").Value = "Betarraga"
$ws.Range("H296").Value = "Sin especificar"
$ws.Range("I296").Value = "Primera"
$ws.Range("J296").Value = 800
$ws.Range("K296").Value = 450
$ws.Range("L296").Value = 500
$ws.Range("M296").Value = 472
$ws.Range("N296").Value = "$/paquete 4 unidades"
$ws.Range("O296").Value = "Región de Arica y Parinacota"
$ws.Range("P296").Value = 118
$ws.Range("Q296").Value = 4
$ws.Range("R296").Value = "Hortaliza"

# --- Row 297 (Segunda / $/paquete 5 unidades) ---
$ws.Range("A297").Value = 1
$ws.Range("B297").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C297").Value = "Arica y Parinacota"
$ws.Range("D297").Value = 44924
$ws.Range("E297").Value = 15
$ws.Range("F297").Value = 100114014
$ws.Range("G297").Value = "Betarraga"
$ws.Range("H297").Value = "Sin especificar"
$ws.Range("I297").Value = "Segunda"
$ws.Range("J297").Value = 600
$ws.Range("K297").Value = 450
$ws.Range("L297").Value = 500
$ws.Range("M297").Value = 471
$ws.Range("N297").Value = "$/paquete 5 unidades"
$ws.Range("O297").Value = "Región de Arica y Parinacota"
$ws.Range("P297").Value = 94
$ws.Range("Q297").Value = 5
$ws.Range("R297").Value = "Hortaliza"
